$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

# Column A holds a literal date-like string ("MM/DD/YYYY") in every other
# row of this sheet (stored as text, not a date serial). Force the cell to
# text formatting before assigning the value so COM's autoconvert doesn't
# turn it into a date serial, then restore the default "Normal" style so
# the cell doesn't pick up a stray explicit style like the rest of the data
# rows.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/06/2026"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 14307.32
$ws.Cells.Item($row, 3).Value = 0.1959103068579821
$ws.Cells.Item($row, 4).Value = 0.8040896931420179
$ws.Cells.Item($row, 5).Value = -54.56
$ws.Cells.Item($row, 6).Value = -9.33
$ws.Cells.Item($row, 7).Value = -19137.39
$ws.Cells.Item($row, 8).Value = -62.46
$ws.Cells.Item($row, 9).Value = -299.9
$ws.Cells.Item($row, 10).Value = -9.67
